$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E/F columns (image canvas width/height) for rows 2-7
$ws.Range("E2:E7").Value = 620
$ws.Range("F2:F7").Value = 586

# Remove the "Done" marker in column I for rows 6, 7, 8, 9
$ws.Range("I6:I9").ClearContents()

# Update selection to I7
$ws.Range("I7").Select()
